# Update column G ("K" - strikeouts) values for rows 2-18 on Sheet1.
# These replace the old pitch-count-derived numbers with the actual
# strikeout totals (regen save_data to use K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 7
    6  = 11
    7  = 6
    8  = 7
    9  = 1
    10 = 5
    11 = 8
    12 = 3
    13 = 6
    14 = 6
    15 = 5
    16 = 5
    17 = 3
    18 = 6
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
